# Apply the "Increase cost of 1-token items to 2 tokens" edit, and fix the
# K/L/M/N formulas on the "Token costs" sheet so that they multiply the
# item's cost (column D) by the presence flag of each hold period, instead
# of the previous (incorrect) "adjacent column" multiplication.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 2: "Token allocation"
# ---------------------------------------------------------------------
$alloc = $wb.Worksheets.Item("Token allocation")
$alloc.Activate()

# The "Ideal" spending example's base-price multiplier was wrong (1
# instead of 2), so both the label and the formula need to be corrected.
$alloc.Range("M18").Value2 = "Base price (2)"
$alloc.Range("N18").Formula = "=N17*2"

$alloc.Range("M21").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 1: "Token costs"
# ---------------------------------------------------------------------
$costs = $wb.Worksheets.Item("Token costs")
$costs.Activate()

# Rows 3-14 hold the individual item costs. Every item that used to cost a
# single token (1) now costs two tokens (2). The 4-token items (Claymore,
# Shield, Health) are unaffected.
$oneTokenRows = @(3, 4, 5, 6, 7, 8, 9, 10, 13)
foreach ($r in $oneTokenRows) {
    $costs.Range("D$r").Value2 = 2
}

# Fix the K/L/M/N formulas for every item/group row (3-16): they should
# multiply the item's cost (column D) by whether it is held during that
# period (columns F, G, H, I), matching column J's D*E pattern - rather
# than the old, buggy "multiply the two neighbouring flag columns" formula.
for ($r = 3; $r -le 16; $r++) {
    $costs.Range("K$r").Formula = "=D$r*F$r"
    $costs.Range("L$r").Formula = "=D$r*G$r"
    $costs.Range("M$r").Formula = "=D$r*H$r"
    $costs.Range("N$r").Formula = "=D$r*I$r"
}

$costs.Range("K11").Select() | Out-Null
